$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.974.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.654.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.509'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.66%  '

$ws.Range("E9").Value = '  +1.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0878'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.887.72'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.662.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.40%  '

$ws.Range("E14").Value = '  +2.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.31'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.982.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0735'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("E20").Value = '  +2.11%  '

$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.28%  '

$ws.Range("E24").Value = '  +1.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.95%  '

$ws.Range("E26").Value = '  +2.21%  '

$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.72%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  +0.46%  '

$ws.Range("E31").Value = '  +2.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.544.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.42%  '

$ws.Range("E35").Value = '  +10.37%  '

$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.582'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.900'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.35%  '

$ws.Range("E39").Value = '  +2.68%  '

$ws.Range("E40").Value = '  +3.64%  '

$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.25%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.952'
$ws.Range("D44").Style = "Normal"

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.774'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.32%  '

$ws.Range("E48").Value = '  +3.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0991'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.29%  '

$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.35%  '
